# Updated cryptos list on Wed Apr  3 17:18:31 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns (and, where the ranked
# coin at a position changed, the Coin (B) / Link (C) columns too) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Column D holds numeric-looking text (e.g. "1.00", "66.284.50", "0.0₃0702").
    # A plain .Value assignment lets Excel auto-coerce that to a Double, which
    # quietly mangles the text (drops trailing zeros, collapses the "." thousands
    # separators, introduces float noise). Force text via NumberFormat "@" while
    # writing, then restore the original style so no stray formatting is left behind.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextCell $ws.Range("D2") '66.284.50'
$ws.Range("E2").Value = '  +1.14%  '
Set-TextCell $ws.Range("D3") '3.346.20'
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextCell $ws.Range("D5") '189.62'
$ws.Range("E5").Value = '  +4.70%  '
Set-TextCell $ws.Range("D6") '558.26'
$ws.Range("E6").Value = '  +0.30%  '
Set-TextCell $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.04%  '
Set-TextCell $ws.Range("D8") '3.341.79'
$ws.Range("E8").Value = '  +2.68%  '
$ws.Range("E9").Value = '  -0.86%  '
Set-TextCell $ws.Range("D10") '0.179'
$ws.Range("E10").Value = '  -3.29%  '
Set-TextCell $ws.Range("D11") '0.582'
$ws.Range("E11").Value = '  -0.53%  '
Set-TextCell $ws.Range("D12") '46.42'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  +1.48%  '
Set-TextCell $ws.Range("D14") '3.882.70'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("E15").Value = '  +0.10%  '
Set-TextCell $ws.Range("D16") '587.87'
$ws.Range("E16").Value = '  -7.62%  '
Set-TextCell $ws.Range("D17") '66.304.08'
$ws.Range("E17").Value = '  +1.11%  '
Set-TextCell $ws.Range("D18") '3.340.80'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws.Range("D19") '0.118'
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws.Range("D20") '17.91'
$ws.Range("E20").Value = '  +1.15%  '
Set-TextCell $ws.Range("D21") '11.01'
$ws.Range("E21").Value = '  -2.94%  '
Set-TextCell $ws.Range("D22") '0.900'
$ws.Range("E22").Value = '  -0.16%  '
Set-TextCell $ws.Range("D23") '18.16'
$ws.Range("E23").Value = '  +2.47%  '
Set-TextCell $ws.Range("D24") '5.03'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +0.01%  '
Set-TextCell $ws.Range("D27") '6.05'
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +2.15%  '
Set-TextCell $ws.Range("D29") '9.48'
$ws.Range("E29").Value = '  -0.43%  '
Set-TextCell $ws.Range("D30") '8.50'
$ws.Range("E30").Value = '  -2.08%  '
Set-TextCell $ws.Range("D31") '30.73'
$ws.Range("E31").Value = '  +1.54%  '
Set-TextCell $ws.Range("D32") '6.67'
$ws.Range("E32").Value = '  +5.83%  '
$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell $ws.Range("D33") '3.82'
$ws.Range("E33").Value = '  -6.08%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws.Range("D34") '583.09'
$ws.Range("E34").Value = '  +5.74%  '
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("E36").Value = '  -0.61%  '
Set-TextCell $ws.Range("D37") '3.776.18'
$ws.Range("E37").Value = '  +4.75%  '
$ws.Range("E38").Value = '  -0.04%  '
Set-TextCell $ws.Range("D39") '56.53'
$ws.Range("E39").Value = '  -0.87%  '
Set-TextCell $ws.Range("D40") '34.38'
$ws.Range("E40").Value = '  +7.50%  '
Set-TextCell $ws.Range("D41") '0.0₃0702'
$ws.Range("E41").Value = '  -2.32%  '
Set-TextCell $ws.Range("D42") '0.127'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("E44").Value = '  -6.84%  '
Set-TextCell $ws.Range("D45") '3.39'
$ws.Range("E45").Value = '  +1.92%  '
Set-TextCell $ws.Range("D46") '0.338'
$ws.Range("E46").Value = '  +0.73%  '
Set-TextCell $ws.Range("D47") '0.0413'
$ws.Range("E47").Value = '  -0.19%  '
Set-TextCell $ws.Range("D48") '2.98'
$ws.Range("E48").Value = '  -17.77%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws.Range("D50") '1.00'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell $ws.Range("D51") '2.55'
$ws.Range("E51").Value = '  -2.10%  '
